# Revert "Layout changes with PinIt routing"
# - Removes the "role" attribute rows (A45:D54) from the Attributes sheet,
#   shifting the remaining rows back up.
# - Fixes the title for p_harris (row 37) back to the "SR. BI Developr" value.
# - Restores the previous sheet selection / scroll position on the Attributes sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")
$ws.Activate()

# Revert the title typo/rename for p_harris (row 37, column C)
$ws.Range("C37").Value = "SR. BI Developr"

# Remove the 10 "role" rows (A45:D54) entirely, shifting rows 55+ up by 10
$ws.Range("A45:D54").EntireRow.Delete()

# Restore prior selection on the sheet
$ws.Range("C55").Select()
